# edit.ps1
# Restructures the "Input" sheet of the PO workbook to match the new
# standard template column layout, and removes a stray empty cell (I3)
# from the 갑지 / 을지 summary sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Input" sheet: rebuild header row + data rows with the new layout
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Input")

# Read the existing data (rows 2-7, columns A..Q) into memory first so
# that overwriting cells in place never clobbers a value we still need
# to read from a different column.
$oldData = @()
for ($r = 2; $r -le 7; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le 17; $c++) {
        $colLetter = [char](64 + $c)
        $rowVals[$colLetter] = $ws.Cells.Item($r, $c).Value2
    }
    $oldData += $rowVals
}

# New header labels, column A..P
$newHeaders = @("발주일자","납기일자","거래처명","거래처 이메일","납품처명","납품처 이메일","프로젝트명","대분류","중분류","소분류","품목명","규격","수량","단가","총금액","비고")

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeaders[$i]
}
# The old header row used a bold / bordered / centered style (s="1");
# the new template header row uses the plain default style.
$ws.Range("A1:P1").ClearFormats()

# Helper: write a text value to a cell while preventing Excel's
# automatic "looks like a date" -> date-serial conversion (needed for
# the order-date / due-date columns, which must stay literal strings).
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

for ($i = 0; $i -lt $oldData.Length; $i++) {
    $r = $i + 2
    $old = $oldData[$i]

    Set-TextCell $ws.Cells.Item($r, 1) $old["C"]   # 발주일자 <- 발주일
    Set-TextCell $ws.Cells.Item($r, 2) $old["D"]   # 납기일자 <- 납기일
    $ws.Cells.Item($r, 3).Value = $old["A"]        # 거래처명 <- 거래처명
    $ws.Cells.Item($r, 4).Value = "영세엔지텍@example.com"  # 거래처 이메일 (new)
    $ws.Cells.Item($r, 5).Value = $old["B"]        # 납품처명 <- 현장명
    $ws.Cells.Item($r, 6).Value = "delivery@example.com"    # 납품처 이메일 (new)
    $ws.Cells.Item($r, 7).Value = $old["B"]        # 프로젝트명 <- 현장명
    $ws.Cells.Item($r, 8).Value = $old["N"]        # 대분류 <- 대분류
    $ws.Cells.Item($r, 9).Value = $old["O"]        # 중분류 <- 중분류
    $ws.Cells.Item($r, 10).Value = $old["P"]       # 소분류 <- 소분류
    $ws.Cells.Item($r, 11).Value = $old["F"]       # 품목명 <- 품목
    $ws.Cells.Item($r, 12).Value = $old["G"]       # 규격 <- 규격
    $ws.Cells.Item($r, 13).Value = $old["H"]       # 수량 <- 수량
    $ws.Cells.Item($r, 14).Value = $old["J"]       # 단가 <- 단가
    $ws.Cells.Item($r, 15).Value = $old["M"]       # 총금액 <- 합계

    $remark = $old["Q"]
    if ($remark -ne $null -and $remark -ne "") {
        $ws.Cells.Item($r, 16).Value = $remark     # 비고 <- 비고
    } else {
        $ws.Cells.Item($r, 16).ClearContents()
    }
}

# Drop the now-unused old column Q entirely so the sheet dimension
# shrinks back down to A1:P7.
$ws.Range("Q1:Q7").ClearContents()
$ws.Range("Q1").EntireColumn.Delete()

# ---------------------------------------------------------------------
# 2. Summary sheets (갑지 / 을지): remove the stray empty cell I3
# ---------------------------------------------------------------------
foreach ($sheetName in @("갑지", "을지")) {
    $sws = $wb.Worksheets.Item($sheetName)
    $sws.Range("I3").ClearContents()
}
